# Add "Norway" and "Poland" market test-data sheets after "Hungary".
# Both new sheets are created the same way the original author built
# them: by duplicating the "Croatia" sheet (which already has the
# right layout/column widths/styles for a single-market panel) and
# then overwriting the two data cells (the NGC ticket code in B4 and
# the market name in B2).

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("Croatia")
$hungary  = $wb.Worksheets.Item("Hungary")

# --- Norway ---------------------------------------------------------
$template.Copy($null, $hungary)
$norway = $wb.Worksheets.Item($wb.Worksheets.Count)
$norway.Name = "Norway"

# Write the code cell before the market-name cell so new shared
# strings land in the same order as the source workbook.
$norway.Range("B4").Value = "NGC-2931/T3601"
$norway.Range("B2").Value = "Norway Market"

# --- Poland ----------------------------------------------------------
$template.Copy($null, $norway)
$poland = $wb.Worksheets.Item($wb.Worksheets.Count)
$poland.Name = "Poland"

$poland.Range("B4").Value = "NGC-2920/T3100/T3104"
$poland.Range("B2").Value = "Poland Market"

# The workbook was left with the "Norway" tab active/selected.
$norway.Activate()
